$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new date values ---
$ws.Range("K1").Value = 20211216
$ws.Range("O1").Value = 20211214

# --- Row 2: header labels (matlab / opencv) ---
$ws.Range("J2").Value = "matlab"
$ws.Range("K2").Value = "matlab"
$ws.Range("L2").Value = "matlab"
$ws.Range("M2").Value = "matlab"
$ws.Range("N2").Value = "opencv"
$ws.Range("O2").Value = "opencv"
$ws.Range("P2").Value = "opencv"
$ws.Range("Q2").Value = "opencv"

# --- Row 3: column labels (fx / fy / cx / cy) ---
$ws.Range("J3").Value = "fx"
$ws.Range("K3").Value = "fy"
$ws.Range("L3").Value = "cx"
$ws.Range("M3").Value = "cy"
$ws.Range("N3").Value = "fx"
$ws.Range("O3").Value = "fy"
$ws.Range("P3").Value = "cx"
$ws.Range("Q3").Value = "cy"

# --- Row 6: new numeric data ---
$ws.Range("J6").Value = 976.55712532041298
$ws.Range("K6").Value = 977.98019880099605
$ws.Range("L6").Value = 673.82457306168806
$ws.Range("M6").Value = 496.58406832349499
$ws.Range("N6").Value = 977.30332833
$ws.Range("O6").Value = 978.96299238999995
$ws.Range("P6").Value = 668.61445031000005
$ws.Range("Q6").Value = 500.89350610000002

# --- Column widths for new bestfit columns K (11) and O (15) ---
$ws.Columns.Item(11).ColumnWidth = 8.714285714285714
$ws.Columns.Item(15).ColumnWidth = 8.714285714285714

# --- Move / resize the chart (shifted down + right) ---
$co = $ws.ChartObjects(1)
$co.Left = 252.5625
$co.Top = 301.5

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("A16").Select() | Out-Null
